$d = $word.ActiveDocument

$replacements = @(
    @{old = "34×17="; new = "61×98="},
    @{old = "91×39="; new = "40×63="},
    @{old = "40×68="; new = "55×67="},
    @{old = "93×27="; new = "16×60="},
    @{old = "47×19="; new = "22×53="},
    @{old = "18×24="; new = "76×82="},
    @{old = "46×53="; new = "16×37="},
    @{old = "68×16="; new = "68×96="},
    @{old = "16×79="; new = "94×64="},
    @{old = "71×41="; new = "73×11="},
    @{old = "61×31="; new = "86×29="},
    @{old = "85×35="; new = "24×30="},
    @{old = "70×23="; new = "42×60="},
    @{old = "20×95="; new = "39×43="},
    @{old = "82×82="; new = "68×42="},
    @{old = "73×34="; new = "35×54="},
    @{old = "38×78="; new = "50×23="},
    @{old = "87×88="; new = "51×11="},
    @{old = "81×76="; new = "59×75="},
    @{old = "41×47="; new = "12×47="},
    @{old = "72×39="; new = "56×35="},
    @{old = "60×48="; new = "98×11="},
    @{old = "60×49="; new = "44×80="},
    @{old = "96×11="; new = "96×18="},
    @{old = "66×25="; new = "55×29="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
